{"js": "// The document contains paragraphs that hold an XML-like \"<id>...</id>\"\n// marker, each split across three runs: \"<id>\", the identifier text\n// (e.g. \"p009r_1\"), and \"</id>\". The edit merges each of these triples\n// into a single run (re-typing the same visible text), which collapses\n// them to one <w:r> carrying the \"<id>\" run's formatting.\n\n// Locate every \"<id>\" opening tag and its matching \"</id>\" closing tag.\nconst openTags = context.document.body.search(\"<id>\", { matchCase: true });\nopenTags.load(\"items\");\nconst closeTags = context.document.body.search(\"</id>\", { matchCase: true });\ncloseTags.load(\"items\");\nawait context.sync();\n\nconst count = Math.min(openTags.items.length, closeTags.items.length);\n\n// Build the combined range for each \"<id>...</id>\" span and load its text.\nconst spans = [];\nfor (let i = 0; i < count; i++) {\n  const span = openTags.items[i].expandTo(closeTags.items[i]);\n  span.load(\"text\");\n  spans.push(span);\n}\nawait context.sync();\n\n// Re-insert the identical text as a replacement; this collapses the\n// three original runs into a single run.\nfor (let i = 0; i < count; i++) {\n  spans[i].insertText(spans[i].text, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The document contains paragraphs that hold an XML-like \"<id>...</id>\"\n# marker, each split across three runs: \"<id>\", the identifier text\n# (e.g. \"p009r_1\"), and \"</id>\". The edit merges each of these triples\n# into a single run (re-typing the same visible text), which collapses\n# them to one run carrying the \"<id>\" run's formatting while keeping the\n# trailing empty run that follows untouched.\n\n$d = $word.ActiveDocument\n\n$cursor = 0\n$iterations = 0\n\nwhile ($true) {\n  # Find the next \"<id>\" opening tag starting at $cursor.\n  $openRng = $d.Content\n  $openRng.Start = $cursor\n  $foundOpen = $openRng.Find.Execute(\"<id>\")\n  if (-not $foundOpen) { break }\n  $idStart = $openRng.Start\n  $idEnd = $openRng.End\n\n  # Find the matching \"</id>\" closing tag right after it.\n  $closeRng = $d.Content\n  $closeRng.Start = $idEnd\n  $foundClose = $closeRng.Find.Execute(\"</id>\")\n  if (-not $foundClose) { break }\n  $closeStart = $closeRng.Start\n  $closeEnd = $closeRng.End\n\n  if ($closeStart -eq $idEnd) {\n    # Already a single span with nothing in between - nothing to merge,\n    # move the cursor past it and keep scanning.\n    $cursor = $closeEnd\n  } else {\n    # Grab the text of the identifier run and the closing-tag run, then\n    # delete both and retype them right after the opening \"<id>\" run.\n    # Retyping collapses the three original runs into a single run that\n    # carries the \"<id>\" run's character formatting.\n    $midText = $d.Range($idEnd, $closeStart).Text\n    $closeText = $d.Range($closeStart, $closeEnd).Text\n\n    $toDelete = $d.Range($idEnd, $closeEnd)\n    $toDelete.Delete()\n\n    $insertPoint = $d.Range($idEnd, $idEnd)\n    $insertPoint.InsertAfter($midText + $closeText)\n\n    $cursor = $idEnd + $midText.Length + $closeText.Length\n  }\n\n  $iterations = $iterations + 1\n  if ($iterations -gt 1000) { break }\n}\n"}
